# Update column F ("dSF") values for the rows that changed during the
# repull/recalculation of the data (see commit message: "repull data,
# push all data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    5  = -4
    8  = 3
    13 = 0
    15 = -1
    18 = -1
    21 = 0
    32 = -6
    33 = 2
    36 = -3
    37 = -7
    38 = -3
    39 = -3
    40 = 0
    42 = 0
    44 = 3
    54 = -1
    55 = 1
    58 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
